$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @('326', '14428723', '2025-08-20', 'Kamil Majchrzak', 'Sebastian Korda', 'Gana Kamil Majchrzak', '2.38'),
    @('327', '14427812', '2025-08-20', 'Marie Bouzkova', 'Beatriz Haddad Maia', 'Gana Beatriz Haddad Maia', '2.75'),
    @('328', '14428006', '2025-08-20', 'Sorana Cirstea', 'Jil Teichmann', 'Gana Jil Teichmann', '3.4'),
    @('329', '14476762', '2025-08-20', 'Arthur Cazaux', 'Jay Clarke', 'Gana Jay Clarke', '5.5'),
    @('330', '14475167', '2025-08-20', 'Cristian Garin', 'Marco Trungelliti', 'Gana Marco Trungelliti', '3.25'),
    @('331', '14475163', '2025-08-20', 'Hugo Grenier', 'Martin Landaluce', 'Gana Hugo Grenier', '3'),
    @('332', '14476745', '2025-08-20', 'Taro Daniel', 'Jan-Lennard Struff', 'Gana Taro Daniel', '3.75'),
    @('333', '14475159', '2025-08-20', 'Matteo Gigante', 'Coleman Wong', 'Gana Matteo Gigante', '2.1'),
    @('334', '14476767', '2025-08-20', 'Patrick Maloney', 'Daniel Elahi Galan', 'Gana Patrick Maloney', '3'),
    @('335', '14477276', '2025-08-20', 'Ugo Blanchet', 'Dmitry Popko', 'Gana Dmitry Popko', '3'),
    @('336', '14485921', '2025-08-20', 'Luca Van Assche', 'Yibing Wu', 'Gana Luca Van Assche', '2.62'),
    @('337', '14477275', '2025-08-20', 'Murphy Cassone', 'Dino Prižmić', 'Gana Murphy Cassone', '3.25'),
    @('338', '14477282', '2025-08-20', 'Fajing Sun', 'Lloyd Harris', 'Gana Fajing Sun', '4'),
    @('339', '14486382', '2025-08-20', 'Juan Manuel Cerundolo', 'Yuta Shimizu', 'Gana Yuta Shimizu', '5.5'),
    @('340', '14485934', '2025-08-20', 'Otto Virtanen', 'Leandro Riedi', 'Gana Otto Virtanen', '2.1'),
    @('341', '14487379', '2025-08-20', 'Alexander Blockx', 'Kimmer Coppejans', 'Gana Kimmer Coppejans', '2.62'),
    @('342', '14475166', '2025-08-20', 'Tessah Andrianjafitrimo', 'Jana Fett', 'Gana Tessah Andrianjafitrimo', '3'),
    @('343', '14476765', '2025-08-20', 'Dominika Salkova', 'Kristina Dmitruk', 'Gana Kristina Dmitruk', '2.25'),
    @('344', '14475168', '2025-08-20', 'Petra Marčinko', 'Alina Charaeva', 'Gana Alina Charaeva', '3.5'),
    @('345', '14477279', '2025-08-20', 'Hina Inoue', 'Simona Waltert', 'Gana Simona Waltert', '1.62'),
    @('346', '14485929', '2025-08-20', 'Viktoriya Tomova', 'Hanne Vandewinkel', 'Gana Hanne Vandewinkel', '3.4'),
    @('347', '14418937', '2025-08-20', 'Stefanos Sakellaridis', 'Maximus Jones', 'Gana Maximus Jones', '2.63')
)

foreach ($row in $rows) {
    $r = $row[0]
    $eventId = $row[1]
    $fecha = $row[2]
    $jugadorA = $row[3]
    $jugadorB = $row[4]
    $pronostico = $row[5]
    $cuota = [double]$row[6]

    # event_id and fecha are written as literal text (leading apostrophe
    # keeps the numeric-looking id and the date string from being
    # auto-converted to a number / date serial), matching the source feed.
    $ws.Range("A$r").Value = "'$eventId"
    $ws.Range("B$r").Value = "'$fecha"
    $ws.Range("C$r").Value = $jugadorA
    $ws.Range("D$r").Value = $jugadorB
    $ws.Range("E$r").Value = $pronostico
    $ws.Range("F$r").Value = $cuota

    # resultado / profit are not known yet for these upcoming matches
    $ws.Range("G$r").Value = "'"
    $ws.Range("H$r").Value = "'"
}
